$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 9246
$ws.Range("I63").Value = 9246
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 9246
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -8622

$ws.Range("H64").Value = 2935.121
$ws.Range("I64").Value = 2652.5
$ws.Range("J64").Value = 3096.6191
$ws.Range("K64").Value = 2652.5
$ws.Range("L64").Value = 3096.6191
$ws.Range("M64").Value = -2404.5
$ws.Range("N64").Value = -3592.6191

$ws.Range("H66").Value = 9246
$ws.Range("I66").Value = 9246
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 27738
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -24618

$ws.Range("H67").Value = 2935.121
$ws.Range("I67").Value = 2652.5
$ws.Range("J67").Value = 3096.6191
$ws.Range("K67").Value = 2652.5
$ws.Range("L67").Value = 3096.6191
$ws.Range("M67").Value = -1794.5
$ws.Range("N67").Value = -4812.6191

$ws.Range("H68").Value = 28222
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 28222
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 28222
$ws.Range("N68").Value = -29720

$ws.Range("H69").Value = 26469
$ws.Range("I69").Value = 20000
$ws.Range("J69").Value = 28086.25
$ws.Range("K69").Value = 60000
$ws.Range("L69").Value = 84258.75
$ws.Range("M69").Value = -59126
$ws.Range("N69").Value = -86006.75

$ws.Range("H70").Value = 1345.75
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1345.75
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4037.25
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -4577.25

$ws.Range("H71").Value = 28222
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 28222
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 84666
$ws.Range("N71").Value = -92154

$ws.Range("H72").Value = 26469
$ws.Range("I72").Value = 20000
$ws.Range("J72").Value = 28086.25
$ws.Range("K72").Value = 180000
$ws.Range("L72").Value = 252776.25
$ws.Range("M72").Value = -175632
$ws.Range("N72").Value = -261512.25

$ws.Range("H73").Value = 1345.75
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1345.75
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4037.25
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -5909.25

$ws.Range("H100").Value = 1971.4286
$ws.Range("I100").Value = 933.3333
$ws.Range("J100").Value = 2750
$ws.Range("K100").Value = 933.3333
$ws.Range("L100").Value = 2750
$ws.Range("M100").Value = -392.3333
$ws.Range("N100").Value = -3832

$ws.Range("H129").Value = 893.21875
$ws.Range("J129").Value = 976.7406999999999
$ws.Range("L129").Value = 2930.2221
$ws.Range("N129").Value = -12930.2221

$ws.Range("H138").Value = 2819164
$ws.Range("I138").Value = 939.8095
$ws.Range("J138").Value = 6900730
$ws.Range("K138").Value = 2819.4285
$ws.Range("L138").Value = 20702190
$ws.Range("M138").Value = 2320.5715
$ws.Range("N138").Value = -20712470

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3494.1
$ws.Range("I32").Value = 2705.5054
$ws.Range("J32").Value = 13971.143
$ws.Range("K32").Value = 2705.5054
$ws.Range("L32").Value = 13971.143
$ws.Range("M32").Value = -2418.5054
$ws.Range("N32").Value = -14545.143

$ws.Range("H118").Value = 30080
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 30080
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 30080
$ws.Range("N118").Value = -33394

$ws.Range("H135").Value = 30000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 30000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H86").Value = 8961.843999999999
$ws.Range("I86").Value = 7888.5454
$ws.Range("J86").Value = 11323.1
$ws.Range("K86").Value = 7888.5454
$ws.Range("L86").Value = 11323.1
$ws.Range("M86").Value = -6765.5454
$ws.Range("N86").Value = -13569.1

$ws.Range("H89").Value = 8961.843999999999
$ws.Range("I89").Value = 7888.5454
$ws.Range("J89").Value = 11323.1
$ws.Range("K89").Value = 39442.727
$ws.Range("L89").Value = 56615.5
$ws.Range("M89").Value = -33826.727
$ws.Range("N89").Value = -67847.5

$ws.Range("H141").Value = 89000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 89000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 89000
$ws.Range("N141").Value = -99360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3800.8333
$ws.Range("I58").Value = 2002.5
$ws.Range("J58").Value = 4700
$ws.Range("K58").Value = 6007.5
$ws.Range("L58").Value = 14100
$ws.Range("M58").Value = -5879.5
$ws.Range("N58").Value = -14356

$ws.Range("H68").Value = 1037.3334
$ws.Range("I68").Value = 667.6667
$ws.Range("J68").Value = 1222.1666
$ws.Range("K68").Value = 2003.0001
$ws.Range("L68").Value = 3666.4998
$ws.Range("M68").Value = -1192.0001
$ws.Range("N68").Value = -5288.4998

$ws.Range("H71").Value = 1037.3334
$ws.Range("I71").Value = 667.6667
$ws.Range("J71").Value = 1222.1666
$ws.Range("K71").Value = 6009.0003
$ws.Range("L71").Value = 10999.4994
$ws.Range("M71").Value = -1953.0003
$ws.Range("N71").Value = -19111.4994

$ws.Range("H113").Value = 1894431
$ws.Range("I113").Value = 2755317
$ws.Range("J113").Value = 481.6
$ws.Range("K113").Value = 8265951
$ws.Range("L113").Value = 1444.8
$ws.Range("M113").Value = -8263781
$ws.Range("N113").Value = -5784.8

$ws.Range("H122").Value = 678.4400000000001
$ws.Range("I122").Value = 569.5714
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 5126.1426
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -2676.1426
$ws.Range("N122").Value = -16150

$ws.Range("H132").Value = 835147.6
$ws.Range("I132").Value = 2158.524
$ws.Range("J132").Value = 2778789
$ws.Range("K132").Value = 19426.716
$ws.Range("L132").Value = 25009101
$ws.Range("M132").Value = -16896.716
$ws.Range("N132").Value = -25014161

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H95").Value = 28755.166
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 28755.166
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 28755.166
$ws.Range("N95").Value = -34247.166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 12999.7
$ws.Range("I68").Value = 22260.4
$ws.Range("J68").Value = 3739
$ws.Range("K68").Value = 22260.4
$ws.Range("L68").Value = 3739
$ws.Range("M68").Value = -21511.4
$ws.Range("N68").Value = -5237

$ws.Range("H71").Value = 12999.7
$ws.Range("I71").Value = 22260.4
$ws.Range("J71").Value = 3739
$ws.Range("K71").Value = 111302
$ws.Range("L71").Value = 18695
$ws.Range("M71").Value = -107558
$ws.Range("N71").Value = -26183

$ws.Range("H132").Value = 3836.9814
$ws.Range("I132").Value = 3729.3157
$ws.Range("K132").Value = 11187.9471
$ws.Range("M132").Value = -8657.947100000001
